$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2218649517684887
$ws.Range("C2").Value = 0.5048231511254019
$ws.Range("J2").Value = 0.03215434083601286
$ws.Range("P2").Value = 0.1543408360128617
$ws.Range("S2").Value = 0.08681672025723473
$ws.Range("B3").Value = 0.01234567901234568
$ws.Range("C3").Value = 0.0308641975308642
$ws.Range("J3").Value = 0.01851851851851852
$ws.Range("P3").Value = 0.7530864197530864
$ws.Range("S3").Value = 0.1851851851851852
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.6222222222222222
$ws.Range("S4").Value = 0.3555555555555556
$ws.Range("B6").Value = 0.05263157894736842
$ws.Range("D6").Value = 0.01435406698564593
$ws.Range("F6").Value = 0.08133971291866028
$ws.Range("J6").Value = 0.2057416267942584
$ws.Range("O6").Value = 0.02392344497607655
$ws.Range("Q6").Value = 0.2200956937799043
$ws.Range("R6").Value = 0.05741626794258373
$ws.Range("S6").Value = 0.3444976076555024
$ws.Range("B7").Value = 0.1096774193548387
$ws.Range("D7").Value = 0.02580645161290323
$ws.Range("F7").Value = 0.05161290322580645
$ws.Range("J7").Value = 0.1032258064516129
$ws.Range("O7").Value = 0.03870967741935484
$ws.Range("Q7").Value = 0.1935483870967742
$ws.Range("R7").Value = 0.07741935483870968
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.09513274336283185
$ws.Range("D8").Value = 0.01769911504424779
$ws.Range("F8").Value = 0.04646017699115045
$ws.Range("J8").Value = 0.1238938053097345
$ws.Range("O8").Value = 0.01769911504424779
$ws.Range("Q8").Value = 0.2057522123893805
$ws.Range("R8").Value = 0.09734513274336283
$ws.Range("S8").Value = 0.3960176991150443
$ws.Range("B9").Value = 0.1333333333333333
$ws.Range("D9").Value = 0.02564102564102564
$ws.Range("E9").Value = 0.005128205128205128
$ws.Range("F9").Value = 0.07179487179487179
$ws.Range("J9").Value = 0.1333333333333333
$ws.Range("O9").Value = 0.005128205128205128
$ws.Range("Q9").Value = 0.1948717948717949
$ws.Range("R9").Value = 0.06153846153846154
$ws.Range("S9").Value = 0.3692307692307693
$ws.Range("B10").Value = 0.105379513633014
$ws.Range("D10").Value = 0.01989683124539425
$ws.Range("F10").Value = 0.06853352984524687
$ws.Range("J10").Value = 0.1171702284450995
$ws.Range("O10").Value = 0.02431834929992631
$ws.Range("Q10").Value = 0.2225497420781135
$ws.Range("R10").Value = 0.0899042004421518
$ws.Range("S10").Value = 0.3522476050110538
$ws.Range("G11").Value = 0.1446280991735537
$ws.Range("J11").Value = 0.07851239669421488
$ws.Range("K11").Value = 0.1776859504132231
$ws.Range("L11").Value = 0.5785123966942148
$ws.Range("S11").Value = 0.02066115702479339
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2285714285714286
$ws.Range("K12").Value = 0.01428571428571429
$ws.Range("L12").Value = 0.02857142857142857
$ws.Range("S12").Value = 0.01428571428571429
$ws.Range("G13").Value = 0.6764705882352942
$ws.Range("J13").Value = 0.2647058823529412
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.02049180327868852
$ws.Range("H15").Value = 0.1680327868852459
$ws.Range("I15").Value = 0.03688524590163934
$ws.Range("J15").Value = 0.3934426229508197
$ws.Range("K15").Value = 0.03278688524590164
$ws.Range("M15").Value = 0.00819672131147541
$ws.Range("O15").Value = 0.06147540983606557
$ws.Range("S15").Value = 0.2786885245901639
$ws.Range("F16").Value = 0.03076923076923077
$ws.Range("H16").Value = 0.1846153846153846
$ws.Range("I16").Value = 0.04615384615384616
$ws.Range("J16").Value = 0.441025641025641
$ws.Range("K16").Value = 0.06153846153846154
$ws.Range("M16").Value = 0.02051282051282051
$ws.Range("O16").Value = 0.08205128205128205
$ws.Range("S16").Value = 0.1333333333333333
$ws.Range("F17").Value = 0.009900990099009901
$ws.Range("H17").Value = 0.1722772277227723
$ws.Range("I17").Value = 0.1049504950495049
$ws.Range("J17").Value = 0.4316831683168317
$ws.Range("K17").Value = 0.09306930693069307
$ws.Range("M17").Value = 0.02376237623762376
$ws.Range("O17").Value = 0.07128712871287128
$ws.Range("S17").Value = 0.09306930693069307
$ws.Range("F18").Value = 0.004975124378109453
$ws.Range("H18").Value = 0.1592039800995025
$ws.Range("I18").Value = 0.06965174129353234
$ws.Range("J18").Value = 0.5174129353233831
$ws.Range("K18").Value = 0.06965174129353234
$ws.Range("M18").Value = 0.02487562189054726
$ws.Range("N18").Value = 0.004975124378109453
$ws.Range("O18").Value = 0.07462686567164178
$ws.Range("S18").Value = 0.07462686567164178
$ws.Range("F19").Value = 0.008842443729903537
$ws.Range("H19").Value = 0.2057877813504823
$ws.Range("I19").Value = 0.08762057877813505
$ws.Range("J19").Value = 0.3987138263665595
$ws.Range("K19").Value = 0.09244372990353698
$ws.Range("M19").Value = 0.01045016077170418
$ws.Range("N19").Value = 0.001607717041800643
$ws.Range("O19").Value = 0.06430868167202572
$ws.Range("S19").Value = 0.1302250803858521
